# "cambios de agosto, puntos fe de ratas e historico"
# Update the reporting period (row 8) from 2021 (Jul-Dec) to 2022 (Jan-Jun),
# refresh the validation dates, and update the saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")
$ws.Activate()

# --- Row 8: reporting period data ---
# A8: Año (year) 2021 -> 2022
$ws.Range("A8").Value = 2022
# B8: fecha de inicio del periodo 2021-07-01 (44378) -> 2022-01-01 (44562)
$ws.Range("B8").Value = 44562
# C8: fecha de termino del periodo 2021-12-31 (44561) -> 2022-06-30 (44742)
$ws.Range("C8").Value = 44742
# V8 / W8: fecha de validacion 2022-01-10 (44571) -> 2022-07-11 (44753)
$ws.Range("V8").Value = 44753
$ws.Range("W8").Value = 44753

# --- Saved view state: scroll position and current selection ---
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 1
$win.ScrollRow = 2
$ws.Range("C10").Select()

# --- Workbook window size (windowWidth 20490 -> 15600 twips, i.e. points*20) ---
$win.Width = 780
$win.Height = 353.25
